$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 134
$ws.Range("H134").Value = 68079.06
$ws.Range("J134").Value = 68079.06
$ws.Range("L134").Value = 68079.06
$ws.Range("N134").Value = -78219.06

# Row 137
$ws.Range("H137").Value = 3196.75
$ws.Range("I137").Value = 2249.8462
$ws.Range("J137").Value = 7300
$ws.Range("K137").Value = 6749.5386
$ws.Range("L137").Value = 21900
$ws.Range("M137").Value = -4199.5386
$ws.Range("N137").Value = -27000

# Row 138
$ws.Range("H138").Value = 2830.75
$ws.Range("I138").Value = 2273.3076
$ws.Range("J138").Value = 3026.6082
$ws.Range("K138").Value = 6819.9228
$ws.Range("L138").Value = 9079.8246
$ws.Range("M138").Value = -1679.9228
$ws.Range("N138").Value = -19359.8246


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 5259
$ws.Range("I32").Value = 4342.5596
$ws.Range("J32").Value = 16256.286
$ws.Range("K32").Value = 4342.5596
$ws.Range("L32").Value = 16256.286
$ws.Range("M32").Value = -4055.5596
$ws.Range("N32").Value = -16830.286

# Row 45
$ws.Range("H45").Value = 16906.7
$ws.Range("I45").Value = 44340.8
$ws.Range("K45").Value = 44340.8
$ws.Range("M45").Value = -43963.8

# Row 102
$ws.Range("H102").Value = 5057.3794
$ws.Range("I102").Value = 4833.231
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 4833.231
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = -3211.231
$ws.Range("N102").Value = -10244

# Row 132
$ws.Range("H132").Value = 2649.9216
$ws.Range("I132").Value = 2309.5134
$ws.Range("K132").Value = 6928.540199999999
$ws.Range("M132").Value = -4398.540199999999


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 3332.5789
$ws.Range("I134").Value = 2893.7273
$ws.Range("K134").Value = 8681.1819
$ws.Range("M134").Value = -6146.1819


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 19
$ws.Range("H19").Value = 630.6
$ws.Range("I19").Value = 645.1111
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 645.1111
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -475.1111
$ws.Range("N19").Value = -840

# Row 24
$ws.Range("H24").Value = 630.6
$ws.Range("I24").Value = 645.1111
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 645.1111
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -475.1111
$ws.Range("N24").Value = -840

# Row 31
$ws.Range("H31").Value = 3498.4854
$ws.Range("I31").Value = 3024.5532
$ws.Range("K31").Value = 3024.5532
$ws.Range("M31").Value = -2729.5532

# Row 34
$ws.Range("H34").Value = 3498.4854
$ws.Range("I34").Value = 3024.5532
$ws.Range("K34").Value = 3024.5532
$ws.Range("M34").Value = -2822.5532

# Row 58
$ws.Range("H58").Value = 3897.4443
$ws.Range("I58").Value = 2612
$ws.Range("K58").Value = 2612
$ws.Range("M58").Value = -2409

# Row 60
$ws.Range("H60").Value = 8087.875
$ws.Range("I60").Value = 5120
$ws.Range("J60").Value = 13034.333
$ws.Range("K60").Value = 5120
$ws.Range("L60").Value = 13034.333
$ws.Range("M60").Value = -4609
$ws.Range("N60").Value = -14056.333

# Row 86
$ws.Range("H86").Value = 2952.4167
$ws.Range("I86").Value = 2270.5557
$ws.Range("K86").Value = 2270.5557
$ws.Range("M86").Value = -1147.5557

# Row 89
$ws.Range("H89").Value = 2952.4167
$ws.Range("I89").Value = 2270.5557
$ws.Range("K89").Value = 11352.7785
$ws.Range("M89").Value = -5736.7785

# Row 136
$ws.Range("H136").Value = 3897.4443
$ws.Range("I136").Value = 2612
$ws.Range("K136").Value = 7836
$ws.Range("M136").Value = -5286

# Row 141
$ws.Range("H141").Value = 346444.88
$ws.Range("J141").Value = 346444.88
$ws.Range("L141").Value = 346444.88
$ws.Range("N141").Value = -356804.88


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 7822.909
$ws.Range("I3").Value = 7765.8887
$ws.Range("K3").Value = 23297.6661
$ws.Range("M3").Value = -23185.6661

# Row 94
$ws.Range("H94").Value = 5390.5
$ws.Range("J94").Value = 6416.5
$ws.Range("L94").Value = 19249.5
$ws.Range("N94").Value = -20601.5

# Row 119
$ws.Range("H119").Value = 729
$ws.Range("I119").Value = 729
$ws.Range("K119").Value = 2187
$ws.Range("M119").Value = 2651

# Row 133
$ws.Range("H133").Value = 20832.666
$ws.Range("I133").Value = 19999.2
$ws.Range("K133").Value = 59997.60000000001
$ws.Range("M133").Value = -54937.60000000001

# Row 134
$ws.Range("H134").Value = 2272
$ws.Range("I134").Value = 1644.091
$ws.Range("J134").Value = 3998.75
$ws.Range("K134").Value = 4932.272999999999
$ws.Range("L134").Value = 11996.25
$ws.Range("M134").Value = 137.7270000000008
$ws.Range("N134").Value = -22136.25

# Row 138
$ws.Range("H138").Value = 8472
$ws.Range("I138").Value = 12500
$ws.Range("J138").Value = 4444
$ws.Range("K138").Value = 37500
$ws.Range("L138").Value = 13332
$ws.Range("M138").Value = -32360
$ws.Range("N138").Value = -23612


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 13
$ws.Range("H13").Value = 249.5
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 122
$ws.Range("H122").Value = 2355.077
$ws.Range("I122").Value = 2426.3333
$ws.Range("K122").Value = 7278.999899999999
$ws.Range("M122").Value = -4828.999899999999


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 6813.6313
$ws.Range("I7").Value = 4825.1
$ws.Range("J7").Value = 9023.111000000001
$ws.Range("K7").Value = 4825.1
$ws.Range("L7").Value = 9023.111000000001
$ws.Range("M7").Value = -4713.1
$ws.Range("N7").Value = -9247.111000000001

# Row 22
$ws.Range("H22").Value = 125001670
$ws.Range("I22").Value = 2092.5
$ws.Range("J22").Value = 500000400
$ws.Range("K22").Value = 2092.5
$ws.Range("L22").Value = 500000400
$ws.Range("M22").Value = -1797.5
$ws.Range("N22").Value = -500000990

# Row 27
$ws.Range("H27").Value = 125001670
$ws.Range("I27").Value = 2092.5
$ws.Range("J27").Value = 500000400
$ws.Range("K27").Value = 2092.5
$ws.Range("L27").Value = 500000400
$ws.Range("M27").Value = -1985.5
$ws.Range("N27").Value = -500000614

# Row 34
$ws.Range("H34").Value = 20000
$ws.Range("I34").Value = 20000
$ws.Range("K34").Value = 20000
$ws.Range("M34").Value = -19828

# Row 61
$ws.Range("H61").Value = 5748.75
$ws.Range("J61").Value = 7331.6665
$ws.Range("L61").Value = 7331.6665
$ws.Range("N61").Value = -7735.6665

# Row 113
$ws.Range("H113").Value = 5748.75
$ws.Range("J113").Value = 7331.6665
$ws.Range("L113").Value = 7331.6665
$ws.Range("N113").Value = -11671.6665

# Row 126
$ws.Range("H126").Value = 6813.6313
$ws.Range("I126").Value = 4825.1
$ws.Range("J126").Value = 9023.111000000001
$ws.Range("K126").Value = 14475.3
$ws.Range("L126").Value = 27069.333
$ws.Range("M126").Value = -12005.3
$ws.Range("N126").Value = -32009.333

# Row 136
$ws.Range("H136").Value = 6679.52
$ws.Range("I136").Value = 5852.8237
$ws.Range("K136").Value = 17558.4711
$ws.Range("M136").Value = -15008.4711


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 54
$ws.Range("H54").Value = 49998
$ws.Range("J54").Value = 49998
$ws.Range("L54").Value = 49998
$ws.Range("N54").Value = -51038

# Row 75
$ws.Range("H75").Value = 49058.5
$ws.Range("I75").Value = 48118
$ws.Range("K75").Value = 48118
$ws.Range("M75").Value = -47182

# Row 78
$ws.Range("H78").Value = 49058.5
$ws.Range("I78").Value = 48118
$ws.Range("K78").Value = 144354
$ws.Range("M78").Value = -139674

# Row 96
$ws.Range("H96").Value = 4776
$ws.Range("I96").Value = 3995.125
$ws.Range("J96").Value = 7899.5
$ws.Range("K96").Value = 3995.125
$ws.Range("L96").Value = 7899.5
$ws.Range("M96").Value = -2622.125
$ws.Range("N96").Value = -10645.5

# Row 107
$ws.Range("H107").Value = 577.2
$ws.Range("I107").Value = 555
$ws.Range("J107").Value = 666
$ws.Range("K107").Value = 1665
$ws.Range("L107").Value = 1998
$ws.Range("M107").Value = 255
$ws.Range("N107").Value = -5838

# Row 132
$ws.Range("H132").Value = 3852.6833
$ws.Range("I132").Value = 3569.7546
$ws.Range("K132").Value = 10709.2638
$ws.Range("M132").Value = -8179.263800000001

# Row 136
$ws.Range("H136").Value = 52633724
$ws.Range("I136").Value = 62501804
$ws.Range("K136").Value = 187505412
$ws.Range("M136").Value = -187502862

